$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "try"
$ws.Range("B18").Value = "try"
$ws.Range("B18").Select()
